# Update the TO-DO List / Completion Tracker sheet:
#  - Row 14 (Wall_Column) gets a Poly Count value and its Measurement Sheet /
#    Model Status columns flip from "NOT STARTED" / "STARTED" to "DONE",
#    matching the styling already used for other completed rows (red bold,
#    centered, no wrap - the same look as cells like F8/F19).
#  - The active selection moves on to the next row (H13).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New poly count for the Wall_Column model.
$ws.Range("D14").Value = 254

# Measurement Sheet status -> DONE
$e14 = $ws.Range("E14")
$e14.Value = "DONE"
$e14.Font.Bold = $true
$e14.Font.Color = 255
$e14.HorizontalAlignment = -4108
$e14.VerticalAlignment = -4108
$e14.WrapText = $false

# Model Status -> DONE (previously "STARTED" with a custom green font/style)
$f14 = $ws.Range("F14")
$f14.Value = "DONE"
$f14.Font.Bold = $true
$f14.Font.Color = 255
$f14.HorizontalAlignment = -4108
$f14.VerticalAlignment = -4108
$f14.WrapText = $false

# Move the selection to H13, matching the saved cursor position.
$ws.Range("H13").Select()
